{"js": "// Update the date header line and the multiplication-table answer grid\n// in place, preserving each run's existing formatting.\n\n// 1) Update the date heading paragraph (first paragraph in the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\n\nif (dateParagraph.text.trim() === \"2024-04-23 Tuesday\") {\n  dateParagraph\n    .getRange()\n    .insertText(\"2024-04-24 Wednesday\", Word.InsertLocation.replace);\n}\n\n// 2) Update the answer cells inside the (only) table. The table has 20\n//    rows x 5 columns; only every 5th row (0, 4, 9, 14, 19) carries the\n//    multiplication answers \u2014 the rows in between are blank spacer rows.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// Old value -> new value, keyed by (row, col) so duplicate text values\n// (e.g. \"92\u00d777=7084\" appears twice, \"27\u00d735=945\" is both a target and a\n// later source) are each replaced independently and correctly.\nconst replacements = {\n  \"0,0\": \"33\u00d717=561\",\n  \"0,1\": \"33\u00d761=2013\",\n  \"0,2\": \"99\u00d742=4158\",\n  \"0,3\": \"25\u00d714=350\",\n  \"0,4\": \"27\u00d735=945\",\n  \"4,0\": \"13\u00d788=1144\",\n  \"4,1\": \"79\u00d774=5846\",\n  \"4,2\": \"68\u00d736=2448\",\n  \"4,3\": \"49\u00d729=1421\",\n  \"4,4\": \"47\u00d747=2209\",\n  \"9,0\": \"72\u00d799=7128\",\n  \"9,1\": \"58\u00d777=4466\",\n  \"9,2\": \"96\u00d774=7104\",\n  \"9,3\": \"91\u00d767=6097\",\n  \"9,4\": \"29\u00d758=1682\",\n  \"14,0\": \"57\u00d728=1596\",\n  \"14,1\": \"62\u00d728=1736\",\n  \"14,2\": \"84\u00d725=2100\",\n  \"14,3\": \"39\u00d790=3510\",\n  \"14,4\": \"97\u00d738=3686\",\n  \"19,0\": \"99\u00d717=1683\",\n  \"19,1\": \"82\u00d784=6888\",\n  \"19,2\": \"86\u00d717=1462\",\n  \"19,3\": \"42\u00d796=4032\",\n  \"19,4\": \"19\u00d736=684\",\n};\n\nconst values = table.values;\nfor (const key of Object.keys(replacements)) {\n  const [r, c] = key.split(\",\").map(Number);\n  if (values[r] && values[r][c] !== undefined) {\n    values[r][c] = replacements[key];\n  }\n}\ntable.values = values;\nawait context.sync();\n", "ps1": "# Update the date header line and the multiplication-table answer grid\n# in place, preserving each run's existing formatting.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date heading paragraph (first paragraph in the body).\n$dateRange = $d.Paragraphs.Item(1).Range\nif ($dateRange.Text.TrimEnd(\"`r`n`a\") -eq \"2024-04-23 Tuesday\") {\n    $dateRange.Text = \"2024-04-24 Wednesday\"\n}\n\n# 2) Update the answer cells inside the (only) table. The table has 20\n#    rows x 5 columns; only every 5th row (1, 5, 10, 15, 20 in 1-based\n#    COM indexing) carries the multiplication answers - the rows in\n#    between are blank spacer rows.\n$t = $d.Tables.Item(1)\n\n# Old value -> new value, keyed by (row, col) so duplicate text values\n# (e.g. \"92x77=7084\" appears twice, \"27x35=945\" is both a target and a\n# later source) are each replaced independently and correctly.\n$newValues = @{\n    \"1,1\"  = \"33\u00d717=561\"\n    \"1,2\"  = \"33\u00d761=2013\"\n    \"1,3\"  = \"99\u00d742=4158\"\n    \"1,4\"  = \"25\u00d714=350\"\n    \"1,5\"  = \"27\u00d735=945\"\n    \"5,1\"  = \"13\u00d788=1144\"\n    \"5,2\"  = \"79\u00d774=5846\"\n    \"5,3\"  = \"68\u00d736=2448\"\n    \"5,4\"  = \"49\u00d729=1421\"\n    \"5,5\"  = \"47\u00d747=2209\"\n    \"10,1\" = \"72\u00d799=7128\"\n    \"10,2\" = \"58\u00d777=4466\"\n    \"10,3\" = \"96\u00d774=7104\"\n    \"10,4\" = \"91\u00d767=6097\"\n    \"10,5\" = \"29\u00d758=1682\"\n    \"15,1\" = \"57\u00d728=1596\"\n    \"15,2\" = \"62\u00d728=1736\"\n    \"15,3\" = \"84\u00d725=2100\"\n    \"15,4\" = \"39\u00d790=3510\"\n    \"15,5\" = \"97\u00d738=3686\"\n    \"20,1\" = \"99\u00d717=1683\"\n    \"20,2\" = \"82\u00d784=6888\"\n    \"20,3\" = \"86\u00d717=1462\"\n    \"20,4\" = \"42\u00d796=4032\"\n    \"20,5\" = \"19\u00d736=684\"\n}\n\nforeach ($rowNum in @(1, 5, 10, 15, 20)) {\n    for ($col = 1; $col -le 5; $col++) {\n        $key = \"$rowNum,$col\"\n        $cell = $t.Cell($rowNum, $col)\n        $cell.Range.Text = $newValues[$key]\n    }\n}\n"}
